$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    "123 Maple Street", "Beverly Hills", "Emma", "Johnson", "+1-310-555-0199", "California", "90210"
)

for ($r = 5; $r -le 7; $r++) {
    for ($c = 1; $c -le 7; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($c -eq 7) {
            $cell.NumberFormat = "@"
            $cell.Value = $data[$c - 1]
            $cell.Style = "Normal"
        } else {
            $cell.Value = $data[$c - 1]
        }
    }
}
